# Generate Report for Archive
#
# The localization status moved from "Ready for handoff" to
# "In Translation" everywhere it is reported:
#   - Overview sheet: per-locale status columns "zh-cn" (E2) and "de-de" (F2)
#   - zh-cn sheet: "Status" column (C2)
#   - de-de sheet: "Status" column (C2)
#
# Because the new text is shorter than the old text, the status columns
# are re-sized to fit the new content (mirrors what Excel's column
# AutoFit does when the cell text changes).

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E1").ColumnWidth = 16.333333333333332
$overview.Range("F1").ColumnWidth = 16.333333333333332

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C1").ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("C1").ColumnWidth = 12.5
